# Update the "S20 Table" worksheet with revised reconstruction-tool
# frequency/percentage data (revision update).

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("S20 Table")

# --- Update the number format used by column C (percentage values) ---
# Previously built-in format "0.00" (numFmtId 2); now a custom "0.0" format.
$ws.Range("C4:C25").NumberFormat = "0.0"

# --- Row data: Metabolic process (A), Frequency (B), Percentage (C) ---
# Rows 11-13 were re-sorted by the updated frequency counts, so the
# "Metabolic process" labels for those three rows are swapped around.

$data = @(
    @{ Row = 4;  Name = "Amino acid metabolism";                          Freq = 188; Pct = 23.095823095823096 }
    @{ Row = 5;  Name = "Carbohydrate metabolism";                        Freq = 156; Pct = 19.164619164619165 }
    @{ Row = 6;  Name = "Metabolism of cofactors and vitamins";           Freq = 88;  Pct = 10.810810810810811 }
    @{ Row = 7;  Name = "Nucleotide metabolism";                          Freq = 60;  Pct = 7.3710073710073711 }
    @{ Row = 8;  Name = "Metabolism of other amino acids";                Freq = 59;  Pct = 7.2481572481572485 }
    @{ Row = 9;  Name = "Energy metabolism";                              Freq = 53;  Pct = 6.5110565110565108 }
    @{ Row = 10; Name = "Lipid metabolism";                               Freq = 37;  Pct = 4.5454545454545459 }
    @{ Row = 11; Name = "Xenobiotics biodegradation and metabolism";      Freq = 29;  Pct = 3.5626535626535625 }
    @{ Row = 12; Name = "ABC transporters";                               Freq = 27;  Pct = 3.3169533169533167 }
    @{ Row = 13; Name = "Biosynthesis of other secondary metabolites";    Freq = 26;  Pct = 3.1941031941031941 }
    @{ Row = 14; Name = "Aminoacyl-tRNA biosynthesis";                    Freq = 20;  Pct = 2.4570024570024569 }
    @{ Row = 15; Name = "Glycan biosynthesis and metabolism";             Freq = 19;  Pct = 2.3341523341523343 }
    @{ Row = 16; Name = "Two-component system";                          Freq = 13;  Pct = 1.597051597051597 }
    @{ Row = 17; Name = "Metabolism of terpenoids and polyketides";       Freq = 13;  Pct = 1.597051597051597 }
    @{ Row = 18; Name = "Vancomycin resistance";                         Freq = 7;   Pct = 0.85995085995085996 }
    @{ Row = 19; Name = "Quorum sensing";                                 Freq = 5;   Pct = 0.61425061425061422 }
    @{ Row = 20; Name = "Sulfur relay system";                           Freq = 5;   Pct = 0.61425061425061422 }
    @{ Row = 21; Name = "Phosphotransferase system (PTS)";               Freq = 4;   Pct = 0.49140049140049141 }
    @{ Row = 22; Name = "Pertussis";                                      Freq = 2;   Pct = 0.24570024570024571 }
    @{ Row = 23; Name = "Bacterial chemotaxis";                           Freq = 1;   Pct = 0.12285012285012285 }
    @{ Row = 24; Name = "Cationic antimicrobial peptide (CAMP) resistance"; Freq = 1; Pct = 0.12285012285012285 }
    @{ Row = 25; Name = "Bacterial secretion system";                     Freq = 1;   Pct = 0.12285012285012285 }
)

foreach ($entry in $data) {
    $r = $entry.Row
    $ws.Cells.Item($r, 1).Value = $entry.Name
    $ws.Cells.Item($r, 2).Value = $entry.Freq
    $ws.Cells.Item($r, 3).Value = $entry.Pct
}

# --- Restore view state: scrolled down with C4:C25 selected ---
$ws.Range("C4:C25").Select()
